$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos crudos")
$ws2 = $wb.Worksheets.Item("Datos válidos")

# New raw timestamps (column C, stored as text) and temperature readings (column E)
# for rows 2..39 of "Datos crudos". Replaces the previous 34 rows (2..35) of data
# with 38 new rows.
$timestamps = @("2023-12-09 06:22:00","2023-12-09 06:23:00","2023-12-09 06:24:00","2023-12-09 06:25:00","2023-12-09 06:26:01","2023-12-09 06:27:01","2023-12-09 06:28:01","2023-12-09 06:29:01","2023-12-09 06:30:02","2023-12-09 06:31:02","2023-12-09 06:32:02","2023-12-09 06:33:02","2023-12-09 06:34:03","2023-12-09 06:35:03","2023-12-09 06:36:03","2023-12-09 06:37:04","2023-12-09 06:38:04","2023-12-09 06:39:04","2023-12-09 06:40:04","2023-12-09 06:41:05","2023-12-09 06:42:05","2023-12-09 06:43:05","2023-12-09 06:44:05","2023-12-09 06:45:06","2023-12-09 06:46:06","2023-12-09 06:47:06","2023-12-09 06:48:07","2023-12-09 06:49:07","2023-12-09 06:50:07","2023-12-09 06:51:07","2023-12-09 06:52:08","2023-12-09 06:53:08","2023-12-09 06:54:08","2023-12-09 06:55:08","2023-12-09 06:56:09","2023-12-09 06:57:09","2023-12-09 06:58:09","2023-12-09 06:59:09")
$evalues = @(25.673076923076898,24.361888111888099,23.444055944055901,23.0506993006993,22.3951048951049,21.870629370629299,21.6083916083915,21.477272727272599,21.346153846153801,20.9527972027971,20.9527972027971,20.5594405594405,20.690559440559401,20.428321678321598,20.428321678321598,20.297202797202701,20.166083916083799,20.166083916083799,20.166083916083799,20.034965034965001,20.166083916083799,20.166083916083799,20.034965034965001,20.166083916083799,20.166083916083799,19.9038461538461,20.034965034965001,19.772727272727199,20.166083916083799,20.034965034965001,20.034965034965001,20.034965034965001,20.034965034965001,20.034965034965001,20.034965034965001,20.034965034965001,20.034965034965001,19.379370629370602)

# Column C holds the raw timestamps as text, so force a text number format
# before writing so Excel does not auto-coerce the strings into date serials.
$ws.Range("C2:C39").NumberFormat = "@"

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    # Rows 36-39 are brand new; fill in the supporting columns (A, B, D) that
    # mirror the pattern used by the rest of the raw-data table.
    if ($row -gt 35) {
        $ws.Cells.Item($row, 1).Value = 6
        $ws.Cells.Item($row, 2).Value = 17
        $ws.Cells.Item($row, 4).Value = 0
    }
    $ws.Cells.Item($row, 3).Value = $timestamps[$i]
    $ws.Cells.Item($row, 5).Value = $evalues[$i]
}

# Manual "fila final" pointer (H2) moved from 16 to 18, and the row-count
# helper formula in H3 no longer subtracts the 2 header offset rows.
$ws.Range("H2").Value = 18
$ws.Range("H3").Formula = "=COUNT(E:E)"

# Chart 1 (on "Datos crudos") plots 'Datos crudos'!$E$2:$E$39 now instead of
# $E$2:$E$35; widening the source range via the Values property keeps the
# c:f reference and cached points in sync automatically. The value axis'
# fixed minimum also moves from 17 to 19.
$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(1).Values = $ws.Range("E2:E39")
$chart1.Axes(2).MinimumScale = 19

# "Datos válidos" sheet: the "T° real max"/"T° real min" helper formulas
# change their offsets from the nominal setpoint (E14).
$ws2.Range("E17").Formula = "=E14+0.7"
$ws2.Range("E20").Formula = "=E14-0.7"
